# daily auto push: 2026-01-27 18:55 UTC
# Insert a new log row for 2026/01/27 (23:00 slot) right before the existing
# row 709, shifting all subsequent rows (old 709..750) down by one
# (new 710..751). This grows the sheet from A1:D750 to A1:D751.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push row 709 (and everything below it) down by one row.
$ws.Rows.Item(709).Insert()

# Fill the newly-inserted row. Column A holds a date written as plain text
# (e.g. "2026/01/27"), so force a text format first to stop Excel from
# auto-converting the string into a date serial, then clear the formatting
# back to the sheet's default (unstyled) cell so the new row matches the
# look of every other data row.
$ws.Range("A709").NumberFormat = "@"
$ws.Range("A709").Value = "2026/01/27"
$ws.Range("A709").ClearFormats()

$ws.Range("B709").Value = "火"
$ws.Range("C709").Value = 23
$ws.Range("D709").Value = 28
